$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@(6, 14, 5, 6)
    ,@(5, 5, 4, 15)
    ,@(4, 6, 3, 14)
    ,@(3, 6, 4, 14)
    ,@(3, 20, 4, 0)
    ,@(4, 8, 5, 12)
    ,@(4, 8, 3, 12)
    ,@(2, 19, 4, 1)
    ,@(8, 6, 2, 14)
    ,@(3, 17, 4, 3)
    ,@(3, 13, 4, 7)
    ,@(2, 2, 3, 18)
    ,@(3, 15, 4, 5)
    ,@(5, 6, 4, 14)
    ,@(4, 3, 3, 17)
    ,@(6, 8, 5, 12)
    ,@(5, 6, 4, 14)
    ,@(3, 15, 4, 5)
    ,@(5, 7, 6, 13)
    ,@(4, 15, 5, 5)
    ,@(2, 6, 3, 14)
    ,@(5, 7, 4, 13)
    ,@(3, 15, 2, 5)
    ,@(4, 6, 3, 14)
    ,@(6, 8, 7, 12)
    ,@(5, 13, 4, 7)
    ,@(4, 5, 3, 15)
    ,@(2, 13, 3, 7)
    ,@(3, 12, 1, 8)
    ,@(4, 12, 6, 8)
    ,@(4, 14, 3, 6)
    ,@(4, 14, 5, 6)
    ,@(3, 4, 5, 16)
    ,@(9, 8, 4, 12)
    ,@(2, 13, 4, 7)
    ,@(5, 14, 6, 6)
    ,@(2, 4, 4, 16)
    ,@(4, 8, 5, 12)
    ,@(2, 13, 3, 7)
    ,@(8, 6, 6, 14)
    ,@(5, 17, 6, 3)
    ,@(4, 16, 3, 4)
    ,@(6, 13, 4, 7)
    ,@(3, 18, 4, 2)
    ,@(4, 4, 2, 16)
    ,@(7, 8, 2, 12)
    ,@(4, 5, 3, 15)
    ,@(4, 1, 3, 19)
    ,@(4, 6, 2, 14)
    ,@(4, 7, 5, 13)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = 1224 + $i
    $rowData = $data[$i]
    $ws.Cells.Item($row, 1).Value = $rowData[0]
    $ws.Cells.Item($row, 2).Value = $rowData[1]
    $ws.Cells.Item($row, 3).Value = $rowData[2]
    $ws.Cells.Item($row, 4).Value = $rowData[3]
}

$ws.Application.ActiveWindow.ScrollRow = 1265
$ws.Range("E1273").Select()
